$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.457784663051898
$ws.Range("C2").Value = 0.814872192099148
$ws.Range("D2").Value = 0.0766847405112316
$ws.Range("E2").Value = 0.410534469403563
$ws.Range("F2").Value = 0.0890782339271882
$ws.Range("G2").Value = 0.619674670797831
$ws.Range("H2").Value = 0.010844306738962
$ws.Range("I2").Value = 0.587141750580945
$ws.Range("J2").Value = 0.922540666150271
$ws.Range("K2").Value = 0.125484120836561
$ws.Range("L2").Value = 0.68009295120062
$ws.Range("M2").Value = 0.419829589465531
$ws.Range("N2").Value = 0.0209140201394268
$ws.Range("O2").Value = 0.103020914020139
$ws.Range("P2").Value = 0.00154918667699458
$ws.Range("Q2").Value = 0.993803253292022
$ws.Range("R2").Value = 0.65143299767622
$ws.Range("S2").Value = 0.00232378001549187
$ws.Range("T2").Value = 0.654531371030209
$ws.Range("U2").Value = 0.635941130906274
$ws.Range("V2").Value = 0.562354763749032
$ws.Range("W2").Value = 0.0449264136328428
$ws.Range("X2").Value = 0.109217660728118
$ws.Range("B3").Value = 0.0565453137103021
$ws.Range("C3").Value = 0.0960495739736638
$ws.Range("D3").Value = 0.341595662277304
$ws.Range("E3").Value = 0.0069713400464756
$ws.Range("F3").Value = 0.297443841982959
$ws.Range("G3").Value = 0.286599535243997
$ws.Range("H3").Value = 0.304415182029435
$ws.Range("I3").Value = 0.0735863671572424
$ws.Range("J3").Value = 0.00929512006196747
$ws.Range("K3").Value = 0.024012393493416
$ws.Range("L3").Value = 0.0286599535243997
$ws.Range("M3").Value = 0.0131680867544539
$ws.Range("N3").Value = 0.0116189000774593
$ws.Range("O3").Value = 0.894655305964369
$ws.Range("P3").Value = 0.92563903950426
$ws.Range("Q3").Value = 0.00464756003098373
$ws.Range("R3").Value = 0.0503485670023238
$ws.Range("S3").Value = 0.477924089852827
$ws.Range("T3").Value = 0.00387296669248644
$ws.Range("U3").Value = 0.0379550735863672
$ws.Range("V3").Value = 0.0658404337722696
$ws.Range("W3").Value = 0.0232378001549187
$ws.Range("X3").Value = 0.0193648334624322
$ws.Range("B4").Value = 0.201394268009295
$ws.Range("C4").Value = 0.0658404337722696
$ws.Range("D4").Value = 0.49419054996127
$ws.Range("E4").Value = 0.554608830364059
$ws.Range("F4").Value = 0.0410534469403563
$ws.Range("G4").Value = 0.0224632068164214
$ws.Range("H4").Value = 0.649109217660728
$ws.Range("I4").Value = 0.0565453137103021
$ws.Range("J4").Value = 0.0511231603408211
$ws.Range("K4").Value = 0.557707203718048
$ws.Range("L4").Value = 0.281177381874516
$ws.Range("M4").Value = 0.556158017041053
$ws.Range("N4").Value = 0.960495739736638
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0542215336948102
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.292021688613478
$ws.Range("S4").Value = 0.010844306738962
$ws.Range("T4").Value = 0.31990704879938
$ws.Range("U4").Value = 0.0255615801704105
$ws.Range("V4").Value = 0.320681642137878
$ws.Range("W4").Value = 0.850503485670023
$ws.Range("X4").Value = 0.817970565453137
$ws.Range("B5").Value = 0.283501161890008
$ws.Range("C5").Value = 0.0224632068164214
$ws.Range("D5").Value = 0.0875290472501936
$ws.Range("E5").Value = 0.0278853601859024
$ws.Range("F5").Value = 0.571649883810999
$ws.Range("G5").Value = 0.0712625871417506
$ws.Range("H5").Value = 0.0356312935708753
$ws.Range("I5").Value = 0.28272656855151
$ws.Range("J5").Value = 0.0170410534469404
$ws.Range("K5").Value = 0.292796281951975
$ws.Range("L5").Value = 0.00929512006196747
$ws.Range("M5").Value = 0.010844306738962
$ws.Range("N5").Value = 0.0069713400464756
$ws.Range("O5").Value = 0.00232378001549187
$ws.Range("P5").Value = 0.0185902401239349
$ws.Range("Q5").Value = 0.00154918667699458
$ws.Range("R5").Value = 0.00542215336948102
$ws.Range("S5").Value = 0.508907823392719
$ws.Range("T5").Value = 0.0216886134779241
$ws.Range("U5").Value = 0.299767621998451
$ws.Range("V5").Value = 0.0511231603408211
$ws.Range("W5").Value = 0.0813323005422153
$ws.Range("X5").Value = 0.0526723470178156
